$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("B2").Value = 14.84080812019566
$ws.Range("C2").Value = 11.26986716156531
$ws.Range("D2").Value = 3.869677050508472
$ws.Range("F2").Value = 15.64784758971977
$ws.Range("G2").Value = 13.76652313815327
$ws.Range("H2").Value = 10.65895552034227
$ws.Range("I2").Value = 14.33961575045192
$ws.Range("O2").Value = 14.19738615010222

# Row 3
$ws.Range("B3").Value = 13.95882645642748
$ws.Range("C3").Value = 10.70066776958554
$ws.Range("D3").Value = 3.745893094442527
$ws.Range("F3").Value = 15.75193083579618
$ws.Range("G3").Value = 13.94792005588849
$ws.Range("H3").Value = 10.72590190416732
$ws.Range("I3").Value = 14.48393346274043
$ws.Range("O3").Value = 14.32336227617294

# Row 4
$ws.Range("B4").Value = 13.38605718177017
$ws.Range("C4").Value = 10.33427683483231
$ws.Range("D4").Value = 3.667924757171807
$ws.Range("F4").Value = 15.82293656527337
$ws.Range("G4").Value = 14.0699823276107
$ws.Range("H4").Value = 10.76937522209283
$ws.Range("I4").Value = 14.57706995224686
$ws.Range("O4").Value = 14.40564674620291

# Row 5
$ws.Range("B5").Value = 13.14488168748407
$ws.Range("C5").Value = 10.18085086907751
$ws.Range("D5").Value = 3.635702219633206
$ws.Range("F5").Value = 15.85364382324834
$ws.Range("G5").Value = 14.12237318977746
$ws.Range("H5").Value = 10.78768686516834
$ws.Range("I5").Value = 14.61616438666758
$ws.Range("O5").Value = 14.440416642005

# Row 6
$ws.Range("B6").Value = 13.10436849939262
$ws.Range("C6").Value = 10.15513015070933
$ws.Range("D6").Value = 3.630325889678673
$ws.Range("F6").Value = 15.858849369712
$ws.Range("G6").Value = 14.13123160592693
$ws.Range("H6").Value = 10.79076350553369
$ws.Range("I6").Value = 14.62272495613541
$ws.Range("O6").Value = 14.4462648756695

# Row 7
$ws.Range("B7").Value = 13.38283592509853
$ws.Range("C7").Value = 10.33222416426993
$ws.Range("D7").Value = 3.667491952858398
$ws.Range("F7").Value = 15.82334353735005
$ws.Range("G7").Value = 14.0706782098413
$ws.Range("H7").Value = 10.76961976584562
$ws.Range("I7").Value = 14.57759257170892
$ws.Range("O7").Value = 14.40611065458507

# Row 8
$ws.Range("B8").Value = 14.54322695243117
$ws.Range("C8").Value = 11.07718628251877
$ws.Range("D8").Value = 3.827429409225175
$ws.Range("F8").Value = 15.68225448076932
$ws.Range("G8").Value = 13.82682889014048
$ws.Range("H8").Value = 10.68154744085763
$ws.Range("I8").Value = 14.38843902702876
$ws.Range("O8").Value = 14.23979697620577

# Row 9
$ws.Range("B9").Value = 16.56823575948581
$ws.Range("C9").Value = 12.39931621284752
$ws.Range("D9").Value = 4.133408836959092
$ws.Range("F9").Value = 15.4624722817201
$ws.Range("G9").Value = 13.43509881480655
$ws.Range("H9").Value = 10.52760511817996
$ws.Range("I9").Value = 14.05328562706742
$ws.Range("O9").Value = 13.95292805383203

# Row 10
$ws.Range("B10").Value = 17.90077219586145
$ws.Range("C10").Value = 13.2811909521447
$ws.Range("D10").Value = 4.358984317764701
$ws.Range("F10").Value = 15.33644467829893
$ws.Range("G10").Value = 13.20234690654872
$ws.Range("H10").Value = 10.42591450149086
$ws.Range("I10").Value = 13.82868073478298
$ws.Range("O10").Value = 13.7662743480966

# Row 11
$ws.Range("B11").Value = 18.47298461884097
$ws.Range("C11").Value = 13.66221254847846
$ws.Range("D11").Value = 4.457204219965171
$ws.Range("F11").Value = 15.28696478313611
$ws.Range("G11").Value = 13.10892836677364
$ws.Range("H11").Value = 10.38212476663279
$ws.Range("I11").Value = 13.73116214270111
$ws.Range("O11").Value = 13.68663537824361

# Row 12
$ws.Range("B12").Value = 18.68476923937947
$ws.Range("C12").Value = 13.80355255330399
$ws.Range("D12").Value = 4.493747108791802
$ws.Range("F12").Value = 15.26936925441486
$ws.Range("G12").Value = 13.07538779268815
$ws.Range("H12").Value = 10.36589759820145
$ws.Range("I12").Value = 13.69490113711167
$ws.Range("O12").Value = 13.65723984207003

# Row 13
$ws.Range("B13").Value = 18.63937572834691
$ws.Range("C13").Value = 13.77324406824677
$ws.Range("D13").Value = 4.48590616190743
$ws.Range("F13").Value = 15.27310780666402
$ws.Range("G13").Value = 13.08252905751286
$ws.Range("H13").Value = 10.36937661939926
$ws.Range("I13").Value = 13.70268095626738
$ws.Range("O13").Value = 13.66353673873425

# Row 14
$ws.Range("B14").Value = 18.49050653153665
$ws.Range("C14").Value = 13.67389992691232
$ws.Range("D14").Value = 4.460223767741364
$ws.Range("F14").Value = 15.2854942423978
$ws.Range("G14").Value = 13.1061319406993
$ws.Range("H14").Value = 10.38078263237279
$ws.Range("I14").Value = 13.72816556961296
$ws.Range("O14").Value = 13.68420169377149

# Row 15
$ws.Range("B15").Value = 18.39868140494679
$ws.Range("C15").Value = 13.61266405782693
$ws.Range("D15").Value = 4.444407286898413
$ws.Range("F15").Value = 15.2932303035256
$ws.Range("G15").Value = 13.12082963782007
$ws.Range("H15").Value = 10.38781538070675
$ws.Range("I15").Value = 13.74386245021611
$ws.Range("O15").Value = 13.69695893683468

# Row 16
$ws.Range("B16").Value = 17.86269180104224
$ws.Range("C16").Value = 13.255880592781
$ws.Range("D16").Value = 4.352475299615676
$ws.Range("F16").Value = 15.33983736227839
$ws.Range("G16").Value = 13.20870646331595
$ws.Range("H16").Value = 10.42882595311042
$ws.Range("I16").Value = 13.83514728631997
$ws.Range("O16").Value = 13.77158534759206

# Row 17
$ws.Range("B17").Value = 17.52516679570333
$ws.Range("C17").Value = 13.03180824639877
$ws.Range("D17").Value = 4.294937973175614
$ws.Range("F17").Value = 15.3704490850473
$ws.Range("G17").Value = 13.26583833488785
$ws.Range("H17").Value = 10.45461710729792
$ws.Range("I17").Value = 13.89233828560832
$ws.Range("O17").Value = 13.81871904028391

# Row 18
$ws.Range("B18").Value = 17.32783481026019
$ws.Range("C18").Value = 12.90103350321513
$ws.Range("D18").Value = 4.261431239014623
$ws.Range("F18").Value = 15.38879436233561
$ws.Range("G18").Value = 13.29986942200534
$ws.Range("H18").Value = 10.4696839735382
$ws.Range("I18").Value = 13.92567134592614
$ws.Range("O18").Value = 13.84632494346568

# Row 19
$ws.Range("B19").Value = 17.2604731796161
$ws.Range("C19").Value = 12.85643169116121
$ws.Range("D19").Value = 4.250016155077716
$ws.Range("F19").Value = 15.39513215977135
$ws.Range("G19").Value = 13.31159146090901
$ws.Range("H19").Value = 10.47482528752023
$ws.Range("I19").Value = 13.93703269737025
$ws.Range("O19").Value = 13.85575686938857

# Row 20
$ws.Range("B20").Value = 17.56142787330138
$ws.Range("C20").Value = 13.0558575532514
$ws.Range("D20").Value = 4.301105786163683
$ws.Range("F20").Value = 15.36711391934495
$ws.Range("G20").Value = 13.25963511157869
$ws.Range("H20").Value = 10.45184753535103
$ws.Range("I20").Value = 13.88620485917544
$ws.Range("O20").Value = 13.81365023066777

# Row 21
$ws.Range("B21").Value = 18.53436609305609
$ws.Range("C21").Value = 13.7031599545737
$ws.Range("D21").Value = 4.467785110092377
$ws.Range("F21").Value = 15.28182497044398
$ws.Range("G21").Value = 13.09914905711375
$ws.Range("H21").Value = 10.37742277445804
$ws.Range("I21").Value = 13.72066203184704
$ws.Range("O21").Value = 13.67811118114257

# Row 22
$ws.Range("B22").Value = 19.14167500361069
$ws.Range("C22").Value = 14.10903064843134
$ws.Range("D22").Value = 4.572919092584408
$ws.Range("F22").Value = 15.23274231583271
$ws.Range("G22").Value = 13.00497787669047
$ws.Range("H22").Value = 10.33085139635326
$ws.Range("I22").Value = 13.6163579896001
$ws.Range("O22").Value = 13.59397221610578

# Row 23
$ws.Range("B23").Value = 18.82015971745848
$ws.Range("C23").Value = 13.89399489862457
$ws.Range("D23").Value = 4.517160395856466
$ws.Range("F23").Value = 15.25832538217487
$ws.Range("G23").Value = 13.0542440948744
$ws.Range("H23").Value = 10.35551808169276
$ws.Range("I23").Value = 13.67167201968627
$ws.Range("O23").Value = 13.63847076839373

# Row 24
$ws.Range("B24").Value = 17.5450444689426
$ws.Range("C24").Value = 13.04499093976803
$ws.Range("D24").Value = 4.298318648837665
$ws.Range("F24").Value = 15.36861942440581
$ws.Range("G24").Value = 13.26243590028396
$ws.Range("H24").Value = 10.45309891376072
$ws.Range("I24").Value = 13.88897636941931
$ws.Range("O24").Value = 13.81594025705461

# Row 25
$ws.Range("B25").Value = 16.04754807621099
$ws.Range("C25").Value = 12.05706855908232
$ws.Range("D25").Value = 4.046581766019021
$ws.Range("F25").Value = 15.51575841838475
$ws.Range("G25").Value = 13.53157622285506
$ws.Range("H25").Value = 10.56724469479352
$ws.Range("I25").Value = 14.14014199453849
$ws.Range("O25").Value = 14.02631273582709
